# "updaite livrable 1 table utilisateur"
#
# The data-dictionary sheet has a "humain" table (rows 29-32, with its
# mirrored/duplicated block in columns H:R) that gets renamed to
# "utilisateur", along with its id column and the hum_* variable-name
# cells becoming uti_*. The "id de l'humain" / "nom de l'humain" / ...
# comment cells in column F keep their original wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 - table name ("humain" -> "utilisateur") and its id column
# ("id_humain" -> "id_utilisateur"), both in the left block (B:F) and the
# mirrored right block (H:R).
$ws.Range("B29").Value = "utilisateur"
$ws.Range("C29").Value = "id_utilisateur"
$ws.Range("H29").Value = "utilisateur"
$ws.Range("I29").Value = "id_utilisateur"

# Row 30 - nom
$ws.Range("C30").Value = "uti_nom"
$ws.Range("I30").Value = "uti_nom"

# Row 31 - prenom
$ws.Range("C31").Value = "uti_prenom"
$ws.Range("I31").Value = "uti_prenom"

# Row 32 - dateNaissance
$ws.Range("C32").Value = "uti_dateNaissance"
$ws.Range("I32").Value = "uti_dateNaissance"

# Leave the sheet with the selection the author ended up with.
$ws.Range("J22").Select()
